$d = $word.ActiveDocument

# Update the date line at the top of the document.
$dateRange = $d.Content
$dateRange.Find.Execute("2025-06-04 Wednesday", $true, $true, $false, $false, $false, $true, 1, $false, "2025-06-05 Thursday", 2)

# Update the multiplication-table cells. Each entry is (row, col, old, new)
# using 1-based Word table addressing; "old" is kept only for documentation/
# traceability with the source diff.
$tbl = $d.Tables.Item(1)

$cellUpdates = @(
    @(1, 1, "435×3=1305", "252×3=756"),
    @(1, 2, "992×4=3968", "820×6=4920"),
    @(1, 3, "839×4=3356", "888×6=5328"),
    @(1, 4, "916×2=1832", "230×4=920"),
    @(1, 5, "329×2=658", "933×5=4665"),
    @(5, 1, "157×4=628", "979×6=5874"),
    @(5, 2, "637×4=2548", "540×9=4860"),
    @(5, 3, "873×8=6984", "416×8=3328"),
    @(5, 4, "423×7=2961", "199×3=597"),
    @(5, 5, "173×6=1038", "734×5=3670"),
    @(10, 1, "261×2=522", "994×2=1988"),
    @(10, 2, "774×3=2322", "566×7=3962"),
    @(10, 3, "403×9=3627", "499×2=998"),
    @(10, 4, "630×7=4410", "610×4=2440"),
    @(10, 5, "298×7=2086", "463×8=3704"),
    @(15, 1, "997×7=6979", "359×9=3231"),
    @(15, 2, "199×3=597", "979×5=4895"),
    @(15, 3, "821×2=1642", "896×8=7168"),
    @(15, 4, "276×6=1656", "927×5=4635"),
    @(15, 5, "414×9=3726", "638×7=4466"),
    @(20, 1, "528×6=3168", "561×7=3927"),
    @(20, 2, "361×6=2166", "578×4=2312"),
    @(20, 3, "412×2=824", "444×5=2220"),
    @(20, 4, "911×7=6377", "674×5=3370"),
    @(20, 5, "609×5=3045", "937×3=2811")
)

foreach ($u in $cellUpdates) {
    $row = $u[0]
    $col = $u[1]
    $new = $u[3]
    # Assign the cell's Range.Text directly (rather than using Find/Replace,
    # which in this runtime replaces every matching occurrence in the whole
    # document instead of just the given Range). Setting Range.Text keeps the
    # run's existing character formatting (font, size, etc.) intact.
    $tbl.Cell($row, $col).Range.Text = $new
}
